$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 with the new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.77 = 36885.53 pesos`n✅ 36885.53 pesos = 8.75 = 942.31 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update tasas sheet values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 114
$wsTasas.Range("O10").Value = 4204.95
$wsTasas.Range("N12").Value = 4214
$wsTasas.Range("O12").Value = 107.655
